$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 33936.53
$ws.Range("J40").Value = 31628.066
$ws.Range("L40").Value = 31628.066
$ws.Range("N40").Value = -31978.066

$ws.Range("H55").Value = 206.75
$ws.Range("I55").Value = 206.75
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 206.75
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 7.25
$ws.Range("N55").Value = ""

$ws.Range("H132").Value = 5406.759
$ws.Range("I132").Value = 1943.8096
$ws.Range("K132").Value = 5831.4288
$ws.Range("M132").Value = -3301.4288

$ws.Range("H138").Value = 6217.511
$ws.Range("J138").Value = 7073.5
$ws.Range("L138").Value = 21220.5
$ws.Range("N138").Value = -31500.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19123.531
$ws.Range("I32").Value = 18541.045
$ws.Range("K32").Value = 18541.045
$ws.Range("M32").Value = -18254.045

$ws.Range("H61").Value = 5031.722
$ws.Range("I61").Value = 3964.3
$ws.Range("K61").Value = 3964.3
$ws.Range("M61").Value = -3752.3

$ws.Range("H74").Value = 32609456
$ws.Range("I74").Value = 41667304
$ws.Range("K74").Value = 41667304
$ws.Range("M74").Value = -41666430

$ws.Range("H77").Value = 32609456
$ws.Range("I77").Value = 41667304
$ws.Range("K77").Value = 208336520
$ws.Range("M77").Value = -208332152

$ws.Range("H97").Value = 596.5833
$ws.Range("I97").Value = 255.44444
$ws.Range("J97").Value = 1620
$ws.Range("K97").Value = 255.44444
$ws.Range("L97").Value = 1620
$ws.Range("M97").Value = 240.55556
$ws.Range("N97").Value = -2612

$ws.Range("H102").Value = 2814.5
$ws.Range("I102").Value = 2814.5
$ws.Range("K102").Value = 2814.5
$ws.Range("M102").Value = -1192.5

$ws.Range("H132").Value = 5489.6943
$ws.Range("I132").Value = 1800.9445
$ws.Range("K132").Value = 5402.833500000001
$ws.Range("M132").Value = -2872.833500000001

$ws.Range("H136").Value = 5031.722
$ws.Range("I136").Value = 3964.3
$ws.Range("K136").Value = 11892.9
$ws.Range("M136").Value = -9342.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 17241928
$ws.Range("I94").Value = 19231152
$ws.Range("K94").Value = 19231152
$ws.Range("M94").Value = -19230701

$ws.Range("H99").Value = 1766.5333
$ws.Range("I99").Value = 1581.3334
$ws.Range("J99").Value = 2507.3333
$ws.Range("K99").Value = 1581.3334
$ws.Range("L99").Value = 2507.3333
$ws.Range("M99").Value = -83.33339999999998
$ws.Range("N99").Value = -5503.3333

$ws.Range("H134").Value = 3576.8823
$ws.Range("I134").Value = 2656.8333
$ws.Range("K134").Value = 7970.499899999999
$ws.Range("M134").Value = -5435.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11906590
$ws.Range("I31").Value = 14085824
$ws.Range("J31").Value = 4620.4614
$ws.Range("K31").Value = 14085824
$ws.Range("L31").Value = 4620.4614
$ws.Range("M31").Value = -14085529
$ws.Range("N31").Value = -5210.4614

$ws.Range("H34").Value = 11906590
$ws.Range("I34").Value = 14085824
$ws.Range("J34").Value = 4620.4614
$ws.Range("K34").Value = 14085824
$ws.Range("L34").Value = 4620.4614
$ws.Range("M34").Value = -14085622
$ws.Range("N34").Value = -5024.4614

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = ""

$ws.Range("H58").Value = 2992
$ws.Range("I58").Value = 2562.2856
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 2562.2856
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -2359.2856
$ws.Range("N58").Value = -6406

$ws.Range("H60").Value = 75000
$ws.Range("J60").Value = 75000
$ws.Range("L60").Value = 75000
$ws.Range("N60").Value = -76022

$ws.Range("H86").Value = 10311.667
$ws.Range("I86").Value = 9769.4
$ws.Range("J86").Value = 10989.5
$ws.Range("K86").Value = 9769.4
$ws.Range("L86").Value = 10989.5
$ws.Range("M86").Value = -8646.4
$ws.Range("N86").Value = -13235.5

$ws.Range("H89").Value = 10311.667
$ws.Range("I89").Value = 9769.4
$ws.Range("J89").Value = 10989.5
$ws.Range("K89").Value = 48847
$ws.Range("L89").Value = 54947.5
$ws.Range("M89").Value = -43231
$ws.Range("N89").Value = -66179.5

$ws.Range("H104").Value = 27999.5
$ws.Range("I104").Value = 25999
$ws.Range("J104").Value = 30000
$ws.Range("K104").Value = 25999
$ws.Range("L104").Value = 30000
$ws.Range("M104").Value = -23378
$ws.Range("N104").Value = -35242

$ws.Range("H121").Value = 59332.668
$ws.Range("J121").Value = 59332.668
$ws.Range("L121").Value = 59332.668
$ws.Range("N121").Value = -61952.668

$ws.Range("H132").Value = 32525200
$ws.Range("I132").Value = 43015360
$ws.Range("J132").Value = 5700.9
$ws.Range("K132").Value = 129046080
$ws.Range("L132").Value = 17102.7
$ws.Range("M132").Value = -129043550
$ws.Range("N132").Value = -22162.7

$ws.Range("H136").Value = 2992
$ws.Range("I136").Value = 2562.2856
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 7686.8568
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -5136.8568
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 180.28572
$ws.Range("I2").Value = 28.666666
$ws.Range("J2").Value = 294
$ws.Range("K2").Value = 171.999996
$ws.Range("L2").Value = 1764
$ws.Range("M2").Value = -58.99999600000001
$ws.Range("N2").Value = -1990

$ws.Range("H10").Value = 276.33334
$ws.Range("I10").Value = 276.33334
$ws.Range("K10").Value = 829.0000200000001
$ws.Range("M10").Value = -690.0000200000001

$ws.Range("H23").Value = 1869
$ws.Range("I23").Value = 1467.75
$ws.Range("J23").Value = 2047.3334
$ws.Range("K23").Value = 4403.25
$ws.Range("L23").Value = 6142.0002
$ws.Range("M23").Value = -4168.25
$ws.Range("N23").Value = -6612.0002

$ws.Range("H113").Value = 1335.3334
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 2997
$ws.Range("M113").Value = -827

$ws.Range("H122").Value = 1600.8
$ws.Range("I122").Value = 1249.5
$ws.Range("J122").Value = 1835
$ws.Range("K122").Value = 11245.5
$ws.Range("L122").Value = 16515
$ws.Range("M122").Value = -8795.5
$ws.Range("N122").Value = -21415

$ws.Range("H141").Value = 7872.2354
$ws.Range("I141").Value = 2989.8572
$ws.Range("K141").Value = 8969.571599999999
$ws.Range("M141").Value = -3789.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 21747492
$ws.Range("I102").Value = 31258210
$ws.Range("K102").Value = 31258210
$ws.Range("M102").Value = -31256588

$ws.Range("H122").Value = 389461
$ws.Range("I122").Value = 1002554
$ws.Range("K122").Value = 3007662
$ws.Range("M122").Value = -3005212

$ws.Range("H132").Value = 3123.8572
$ws.Range("I132").Value = 2798.6572
$ws.Range("J132").Value = 3936.8572
$ws.Range("K132").Value = 8395.971600000001
$ws.Range("L132").Value = 11810.5716
$ws.Range("M132").Value = -5865.971600000001
$ws.Range("N132").Value = -16870.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3788.0864
$ws.Range("I7").Value = 3179.691
$ws.Range("K7").Value = 3179.691
$ws.Range("M7").Value = -3067.691

$ws.Range("H10").Value = 1795.75
$ws.Range("I10").Value = 841.5
$ws.Range("J10").Value = 2750
$ws.Range("K10").Value = 841.5
$ws.Range("L10").Value = 2750
$ws.Range("M10").Value = -701.5
$ws.Range("N10").Value = -3030

$ws.Range("H40").Value = 16371854
$ws.Range("I40").Value = 6252023.5
$ws.Range("J40").Value = 41671428
$ws.Range("K40").Value = 6252023.5
$ws.Range("L40").Value = 41671428
$ws.Range("M40").Value = -6251887.5
$ws.Range("N40").Value = -41671700

$ws.Range("H46").Value = 3554.4827
$ws.Range("I46").Value = 3045.1428
$ws.Range("J46").Value = 3716.5454
$ws.Range("K46").Value = 3045.1428
$ws.Range("L46").Value = 3716.5454
$ws.Range("M46").Value = -2857.1428
$ws.Range("N46").Value = -4092.5454

$ws.Range("H122").Value = 8058.1304
$ws.Range("I122").Value = 3770.7144
$ws.Range("J122").Value = 9933.875
$ws.Range("K122").Value = 11312.1432
$ws.Range("L122").Value = 29801.625
$ws.Range("M122").Value = -8862.143199999999
$ws.Range("N122").Value = -34701.625

$ws.Range("H126").Value = 3788.0864
$ws.Range("I126").Value = 3179.691
$ws.Range("K126").Value = 9539.073
$ws.Range("M126").Value = -7069.073

$ws.Range("H140").Value = 51457.57
$ws.Range("J140").Value = 51457.57
$ws.Range("L140").Value = 51457.57
$ws.Range("N140").Value = -61817.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 62802
$ws.Range("J46").Value = 62802
$ws.Range("L46").Value = 62802
$ws.Range("N46").Value = -63264

$ws.Range("H100").Value = 1025.0667
$ws.Range("I100").Value = 746.4
$ws.Range("J100").Value = 1582.4
$ws.Range("K100").Value = 1492.8
$ws.Range("L100").Value = 3164.8
$ws.Range("M100").Value = -951.8
$ws.Range("N100").Value = -4246.8

$ws.Range("H126").Value = 1185.5294
$ws.Range("I126").Value = 1143.6666
$ws.Range("J126").Value = 1499.5
$ws.Range("K126").Value = 3430.9998
$ws.Range("L126").Value = 4498.5
$ws.Range("M126").Value = -960.9998000000001
$ws.Range("N126").Value = -9438.5

$ws.Range("H132").Value = 3965.5
$ws.Range("I132").Value = 3232.8333
$ws.Range("K132").Value = 9698.499899999999
$ws.Range("M132").Value = -7168.499899999999

$ws.Range("H134").Value = 62802
$ws.Range("J134").Value = 62802
$ws.Range("L134").Value = 188406
$ws.Range("N134").Value = -193476
